$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header/label cells (shared strings: "Average flash", "Average HTML5")
$ws.Range("G2").Value = "Average flash"
$ws.Range("G3").Value = "Average HTML5"

# New average formulas, styled like the existing percentage column (style index 1)
$ws.Range("H2").Formula = "=(E1+E3+E5)/3"
$ws.Range("H3").Formula = "=(E2+E4+E6)/3"

$ws.Range("H2:H3").NumberFormat = "0%"

# Set column G width explicitly (matches added <cols> entry, closest
# representable width to Excel's default 9.140625 given this engine's
# internal pixel-quantised ColumnWidth rounding)
$ws.Columns.Item(7).ColumnWidth = 8.28

# Move the active selection to H3, matching the authored selection change
$ws.Range("H3").Select()
